# Refresh cached market-price-derived figures (currentAveragePrice*, LevePrice*,
# LeveProfit* columns H:N) on the Leve-profit tables for a batch of rows across
# all eight job sheets, as produced by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 480.03845
$ws.Range("I33").Value = 154.71428
$ws.Range("J33").Value = 1846.4
$ws.Range("K33").Value = 154.71428
$ws.Range("L33").Value = 1846.4
$ws.Range("M33").Value = 74.28572
$ws.Range("N33").Value = -2304.4

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 20463.156
$ws.Range("I132").Value = 2619.5112
$ws.Range("J132").Value = 154290.5
$ws.Range("K132").Value = 7858.5336
$ws.Range("L132").Value = 462871.5
$ws.Range("M132").Value = -5328.5336
$ws.Range("N132").Value = -467931.5

# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 50001730
$ws.Range("J135").Value = 250000510
$ws.Range("L135").Value = 2250004590
$ws.Range("N135").Value = -2250009660

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 6931.727
$ws.Range("I137").Value = 10515.692
$ws.Range("J137").Value = 5428.7744
$ws.Range("K137").Value = 31547.076
$ws.Range("L137").Value = 16286.3232
$ws.Range("M137").Value = -28997.076
$ws.Range("N137").Value = -21386.3232


$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 12839.532
$ws.Range("I32").Value = 12728.514
$ws.Range("J32").Value = 13250.3
$ws.Range("K32").Value = 12728.514
$ws.Range("L32").Value = 13250.3
$ws.Range("M32").Value = -12441.514
$ws.Range("N32").Value = -13824.3

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1599.1132
$ws.Range("I74").Value = 1410.9762
$ws.Range("K74").Value = 1410.9762
$ws.Range("M74").Value = -536.9762000000001

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1599.1132
$ws.Range("I77").Value = 1410.9762
$ws.Range("K77").Value = 7054.881
$ws.Range("M77").Value = -2686.881

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 10640218
$ws.Range("I132").Value = 13890580
$ws.Range("J132").Value = 2672
$ws.Range("K132").Value = 41671740
$ws.Range("L132").Value = 8016
$ws.Range("M132").Value = -41669210
$ws.Range("N132").Value = -13076


$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2659.3103
$ws.Range("I134").Value = 2853.875
$ws.Range("J134").Value = 2419.8462
$ws.Range("K134").Value = 8561.625
$ws.Range("L134").Value = 7259.5386
$ws.Range("M134").Value = -6026.625
$ws.Range("N134").Value = -12329.5386


$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 4019549.5
$ws.Range("I31").Value = 1403.0975
$ws.Range("J31").Value = 7942025.5
$ws.Range("K31").Value = 1403.0975
$ws.Range("L31").Value = 7942025.5
$ws.Range("M31").Value = -1108.0975
$ws.Range("N31").Value = -7942615.5

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 4019549.5
$ws.Range("I34").Value = 1403.0975
$ws.Range("J34").Value = 7942025.5
$ws.Range("K34").Value = 1403.0975
$ws.Range("L34").Value = 7942025.5
$ws.Range("M34").Value = -1201.0975
$ws.Range("N34").Value = -7942429.5

# Row 94: Beech, Please / Beech Lumber
$ws.Range("H94").Value = 1427.579
$ws.Range("I94").Value = 547.2
$ws.Range("J94").Value = 1742
$ws.Range("K94").Value = 547.2
$ws.Range("L94").Value = 1742
$ws.Range("M94").Value = -96.20000000000005
$ws.Range("N94").Value = -2644


$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 10944253
$ws.Range("I68").Value = 12346942
$ws.Range("J68").Value = 10418245
$ws.Range("K68").Value = 37040826
$ws.Range("L68").Value = 31254735
$ws.Range("M68").Value = -37040015
$ws.Range("N68").Value = -31256357

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 10944253
$ws.Range("I71").Value = 12346942
$ws.Range("J71").Value = 10418245
$ws.Range("K71").Value = 111122478
$ws.Range("L71").Value = 93764205
$ws.Range("M71").Value = -111118422
$ws.Range("N71").Value = -93772317

# Row 103: West Meats East / Nomad Meat Pie
$ws.Range("H103").Value = 964.8570999999999
$ws.Range("I103").Value = 1500
$ws.Range("J103").Value = 563.5
$ws.Range("K103").Value = 4500
$ws.Range("L103").Value = 1690.5
$ws.Range("M103").Value = -3621
$ws.Range("N103").Value = -3448.5

# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 3949.1562
$ws.Range("I113").Value = 5003.609
$ws.Range("J113").Value = 1254.4445
$ws.Range("K113").Value = 15010.827
$ws.Range("L113").Value = 3763.3335
$ws.Range("M113").Value = -12840.827
$ws.Range("N113").Value = -8103.333500000001

# Row 120: A Happy End / Paella
$ws.Range("H120").Value = 384111.25
$ws.Range("I120").Value = 1501445
$ws.Range("J120").Value = 11666.667
$ws.Range("K120").Value = 4504335
$ws.Range("L120").Value = 35000.001
$ws.Range("M120").Value = -4499497
$ws.Range("N120").Value = -44676.001

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 4543.727
$ws.Range("I131").Value = 17193.166
$ws.Range("J131").Value = 1732.7407
$ws.Range("K131").Value = 51579.49800000001
$ws.Range("L131").Value = 5198.2221
$ws.Range("M131").Value = -46539.49800000001
$ws.Range("N131").Value = -15278.2221


$ws = $wb.Worksheets.Item("GSM")
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 21674
$ws.Range("I126").Value = 80502.39999999999
$ws.Range("J126").Value = 2064.5334
$ws.Range("K126").Value = 241507.2
$ws.Range("L126").Value = 6193.600199999999
$ws.Range("M126").Value = -239037.2
$ws.Range("N126").Value = -11133.6002


$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 58827050
$ws.Range("I7").Value = 90911000
$ws.Range("J7").Value = 6468.3335
$ws.Range("K7").Value = 90911000
$ws.Range("L7").Value = 6468.3335
$ws.Range("M7").Value = -90910888
$ws.Range("N7").Value = -6692.3335

# Row 16: Saddle Sore / Hard Leather
$ws.Range("H16").Value = 2257.9565
$ws.Range("I16").Value = 2351.8
$ws.Range("J16").Value = 1632.3334
$ws.Range("K16").Value = 2351.8
$ws.Range("L16").Value = 1632.3334
$ws.Range("M16").Value = -2181.8
$ws.Range("N16").Value = -1972.3334

# Row 35: No Risk, No Reward / Toadskin Cesti
$ws.Range("H35").Value = 2098.375
$ws.Range("I35").Value = 2098.375
$ws.Range("K35").Value = 2098.375
$ws.Range("M35").Value = -1762.375

# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 4830.7
$ws.Range("I40").Value = 3750.8333
$ws.Range("J40").Value = 6450.5
$ws.Range("K40").Value = 3750.8333
$ws.Range("L40").Value = 6450.5
$ws.Range("M40").Value = -3614.8333
$ws.Range("N40").Value = -6722.5

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 1301.3334
$ws.Range("I61").Value = 1104
$ws.Range("K61").Value = 1104
$ws.Range("M61").Value = -902

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 7579420
$ws.Range("I82").Value = 3564
$ws.Range("J82").Value = 16670447
$ws.Range("K82").Value = 3564
$ws.Range("L82").Value = 16670447
$ws.Range("M82").Value = -3203
$ws.Range("N82").Value = -16671169

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 7579420
$ws.Range("I85").Value = 3564
$ws.Range("J85").Value = 16670447
$ws.Range("K85").Value = 3564
$ws.Range("L85").Value = 16670447
$ws.Range("M85").Value = -2316
$ws.Range("N85").Value = -16672943

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1301.3334
$ws.Range("I113").Value = 1104
$ws.Range("K113").Value = 1104
$ws.Range("M113").Value = 1066

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 69220.266
$ws.Range("I122").Value = 113378.22
$ws.Range("J122").Value = 2983.3333
$ws.Range("K122").Value = 340134.66
$ws.Range("L122").Value = 8949.999899999999
$ws.Range("M122").Value = -337684.66
$ws.Range("N122").Value = -13849.9999

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 58827050
$ws.Range("I126").Value = 90911000
$ws.Range("J126").Value = 6468.3335
$ws.Range("K126").Value = 272733000
$ws.Range("L126").Value = 19405.0005
$ws.Range("M126").Value = -272730530
$ws.Range("N126").Value = -24345.0005


$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1851.841
$ws.Range("I132").Value = 1410.3572
$ws.Range("K132").Value = 4231.0716
$ws.Range("M132").Value = -1701.0716
